$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.2323342713546026
$ws.Range("D2").Value = 0.1936821668427342
$ws.Range("E2").Value = 0.1637286036555636
$ws.Range("F2").Value = 1.371863024247688
$ws.Range("G2").Value = 0.7656368203074919
$ws.Range("H2").Value = 0.8573159392711887
$ws.Range("I2").Value = 0.9486264204944845
$ws.Range("J2").Value = 0.1817389190539274
$ws.Range("K2").Value = 1.239856818260819
$ws.Range("L2").Value = 0.2145138343572981
$ws.Range("N2").Value = 1.237564196572151
$ws.Range("O2").Value = 3.257304884545078

# Row 3
$ws.Range("C3").Value = 0.2288162730677783
$ws.Range("D3").Value = 0.1904205364025984
$ws.Range("E3").Value = 0.1626116673094309
$ws.Range("F3").Value = 1.375586268995065
$ws.Range("G3").Value = 0.7687157380360929
$ws.Range("H3").Value = 0.8628302847912366
$ws.Range("I3").Value = 0.9515894264282423
$ws.Range("J3").Value = 0.1815991764158227
$ws.Range("K3").Value = 1.119980459949659
$ws.Range("L3").Value = 0.2139799278543535
$ws.Range("N3").Value = 1.232029302888833
$ws.Range("O3").Value = 3.275092097228779

# Row 4
$ws.Range("C4").Value = 0.2267406629446072
$ws.Range("D4").Value = 0.1884801879155162
$ws.Range("E4").Value = 0.1619849363818346
$ws.Range("F4").Value = 1.378574404789475
$ws.Range("G4").Value = 0.7710850826054596
$ws.Range("H4").Value = 0.866578389902088
$ws.Range("I4").Value = 0.9539019835516811
$ws.Range("J4").Value = 0.1815829987920772
$ws.Range("K4").Value = 1.046306572079743
$ws.Range("L4").Value = 0.2137279907519272
$ws.Range("N4").Value = 1.229038101946102
$ws.Range("O4").Value = 3.287777044213598

# Row 5
$ws.Range("C5").Value = 0.2259161870615003
$ws.Range("D5").Value = 0.1877052544943609
$ws.Range("E5").Value = 0.1617444482589399
$ws.Range("F5").Value = 1.379968744233828
$ws.Range("G5").Value = 0.7721709803561723
$ws.Range("H5").Value = 0.8681969366675162
$ws.Range("I5").Value = 0.9549684688646352
$ws.Range("J5").Value = 0.1815939521987922
$ws.Range("K5").Value = 1.01626894808922
$ws.Range("L5").Value = 0.2136444608577079
$ws.Range("N5").Value = 1.227922038103159
$ws.Range("O5").Value = 3.293389749473732

# Row 6
$ws.Range("C6").Value = 0.2257805764718768
$ws.Range("D6").Value = 0.1875775332706553
$ws.Range("E6").Value = 0.1617054175399275
$ws.Range("F6").Value = 1.380210945809637
$ws.Range("G6").Value = 0.7723585615662216
$ws.Range("H6").Value = 0.8684712031382773
$ws.Range("I6").Value = 0.9551530551692977
$ws.Range("J6").Value = 0.1815968318573233
$ws.Range("K6").Value = 1.011280396653831
$ws.Range("L6").Value = 0.2136317481979972
$ws.Range("N6").Value = 1.227742943858942
$ws.Range("O6").Value = 3.294348521389097

# Row 7
$ws.Range("C7").Value = 0.2267294571676928
$ws.Range("D7").Value = 0.1884696728800463
$ws.Range("E7").Value = 0.1619816326315799
$ws.Range("F7").Value = 1.37859249393945
$ws.Range("G7").Value = 0.7710992400889154
$ws.Range("H7").Value = 0.8665998489839808
$ws.Range("I7").Value = 0.9539158639849887
$ws.Range("J7").Value = 0.1815830754199297
$ws.Range("K7").Value = 1.04590153123948
$ws.Range("L7").Value = 0.213726786682102
$ws.Range("N7").Value = 1.229022633155608
$ws.Range("O7").Value = 3.287850943417951

# Row 8
$ws.Range("C8").Value = 0.2311038044655618
$ws.Range("D8").Value = 0.1925446879597388
$ws.Range("E8").Value = 0.1633312561705615
$ws.Range("F8").Value = 1.37300115489375
$ws.Range("G8").Value = 0.7665990108612064
$ws.Range("H8").Value = 0.8591421351379296
$ws.Range("I8").Value = 0.9495457265831178
$ws.Range("J8").Value = 0.1816763076152981
$ws.Range("K8").Value = 1.198539307054062
$ws.Range("L8").Value = 0.2143140261324348
$ws.Range("N8").Value = 1.235571516979832
$ws.Range("O8").Value = 3.263071912015192

# Row 9
$ws.Range("C9").Value = 0.2403475519714107
$ws.Range("D9").Value = 0.2010259776695307
$ws.Range("E9").Value = 0.1664444058852439
$ws.Range("F9").Value = 1.367603191523571
$ws.Range("G9").Value = 0.7615765805386303
$ws.Range("H9").Value = 0.8473893396532617
$ws.Range("I9").Value = 0.9448876290246133
$ws.Range("J9").Value = 0.1824100802224748
$ws.Range("K9").Value = 1.497215412687808
$ws.Range("L9").Value = 0.2160655281615931
$ws.Range("N9").Value = 1.251626451801201
$ws.Range("O9").Value = 3.228473706458146

# Row 10
$ws.Range("C10").Value = 0.2475396670640464
$ws.Range("D10").Value = 0.2075511956291365
$ws.Range("E10").Value = 0.1690135308426832
$ws.Range("F10").Value = 1.367027551500229
$ws.Range("G10").Value = 0.7602096753743695
$ws.Range("H10").Value = 0.8405018871650611
$ws.Range("I10").Value = 0.9438484431157832
$ws.Range("J10").Value = 0.1832833039261814
$ws.Range("K10").Value = 1.716150863001474
$ws.Range("L10").Value = 0.2177155167091271
$ws.Range("N10").Value = 1.265358122103805
$ws.Range("O10").Value = 3.211589051451341

# Row 11
$ws.Range("C11").Value = 0.250897541046271
$ws.Range("D11").Value = 0.2105825539028814
$ws.Range("E11").Value = 0.1702429928478395
$ws.Range("F11").Value = 1.367501252957069
$ws.Range("G11").Value = 0.7600933979581157
$ws.Range("H11").Value = 0.8377473388483736
$ws.Range("I11").Value = 0.9438929098891649
$ws.Range("J11").Value = 0.1837527540491948
$ws.Range("K11").Value = 1.81561971423929
$ws.Range("L11").Value = 0.2185444738616766
$ws.Range("N11").Value = 1.272020959701123
$ws.Range("O11").Value = 3.205762058807608

# Row 12
$ws.Range("C12").Value = 0.2521813631977921
$ws.Range("D12").Value = 0.2117394119023004
$ws.Range("E12").Value = 0.1707172437348241
$ws.Range("F12").Value = 1.367786324069755
$ws.Range("G12").Value = 0.760122134759996
$ws.Range("H12").Value = 0.8367586469161807
$ws.Range("I12").Value = 0.9439840787447125
$ws.Range("J12").Value = 0.1839408728666498
$ws.Range("K12").Value = 1.853265586113082
$ws.Range("L12").Value = 0.2188695973972727
$ws.Range("N12").Value = 1.274603446909268
$ws.Range("O12").Value = 3.203822162535545

# Row 13
$ws.Range("C13").Value = 0.2519043253985984
$ws.Range("D13").Value = 0.211489865624273
$ws.Range("E13").Value = 0.1706147201545711
$ws.Range("F13").Value = 1.367720229607215
$ws.Range("G13").Value = 0.7601127082207029
$ws.Range("H13").Value = 0.8369691608610452
$ws.Range("I13").Value = 0.9439611387206739
$ws.Range("J13").Value = 0.1838998984512941
$ws.Range("K13").Value = 1.845158843184322
$ws.Range("L13").Value = 0.218799078274948
$ws.Range("N13").Value = 1.27404462648596
$ws.Range("O13").Value = 3.204228093104604

# Row 14
$ws.Range("C14").Value = 0.2510029166790417
$ws.Range("D14").Value = 0.210677550605709
$ws.Range("E14").Value = 0.1702818361084795
$ws.Range("F14").Value = 1.36752258820627
$ws.Range("G14").Value = 0.76009430334652
$ws.Range("H14").Value = 0.8376649086531245
$ws.Range("I14").Value = 0.9438989209605495
$ws.Range("J14").Value = 0.1837680235116608
$ws.Range("K14").Value = 1.818717295586794
$ws.Range("L14").Value = 0.2185709975051608
$ws.Range("N14").Value = 1.272232234196494
$ws.Range("O14").Value = 3.205597117742485

# Row 15
$ws.Range("C15").Value = 0.2504523719443199
$ws.Range("D15").Value = 0.2101811462412542
$ws.Range("E15").Value = 0.1700790637871705
$ws.Range("F15").Value = 1.367415288722569
$ws.Range("G15").Value = 0.7600925085961023
$ws.Range("H15").Value = 0.8380981568674883
$ws.Range("I15").Value = 0.9438704895170673
$ws.Range("J15").Value = 0.1836885928709577
$ws.Range("K15").Value = 1.802518292165018
$ws.Range("L15").Value = 0.2184327504312762
$ws.Range("N15").Value = 1.271129816302491
$ws.Range("O15").Value = 3.206470413758268

# Row 16
$ws.Range("C16").Value = 0.2473219453514162
$ws.Range("D16").Value = 0.207354347506481
$ws.Range("E16").Value = 0.1689344003454316
$ws.Range("F16").Value = 1.36701138838005
$ws.Range("G16").Value = 0.7602274533153377
$ws.Range("H16").Value = 0.8406895177075739
$ws.Range("I16").Value = 0.9438559408594003
$ws.Range("J16").Value = 0.1832540743234929
$ws.Range("K16").Value = 1.709647572963206
$ws.Range("L16").Value = 0.2176629147855849
$ws.Range("N16").Value = 1.264931028018353
$ws.Range("O16").Value = 3.212007170604238

# Row 17
$ws.Range("C17").Value = 0.2454235088372627
$ws.Range("D17").Value = 0.2056362555076277
$ws.Range("E17").Value = 0.1682477071020401
$ws.Range("F17").Value = 1.366951945581818
$ws.Range("G17").Value = 0.7604397695672134
$ws.Range("H17").Value = 0.8423761672801646
$ws.Range("I17").Value = 0.9439794526806509
$ws.Range("J17").Value = 0.1830059817062875
$ws.Range("K17").Value = 1.65264019905112
$ws.Range("L17").Value = 0.2172106793711848
$ws.Range("N17").Value = 1.261234565109319
$ws.Range("O17").Value = 3.215878674254782

# Row 18
$ws.Range("C18").Value = 0.2443396938528508
$ws.Range("D18").Value = 0.2046539913167607
$ws.Range("E18").Value = 0.1678584611329015
$ws.Range("F18").Value = 1.366986995308096
$ws.Range("G18").Value = 0.7606094672578649
$ws.Range("H18").Value = 0.8433819206992155
$ws.Range("I18").Value = 0.9440991770632223
$ws.Range("J18").Value = 0.1828700849066962
$ws.Range("K18").Value = 1.619839368581893
$ws.Range("L18").Value = 0.2169579452133519
$ws.Range("N18").Value = 1.259147660911736
$ws.Range("O18").Value = 3.2182799567112

# Row 19
$ws.Range("C19").Value = 0.2439741302018348
$ws.Range("D19").Value = 0.2043224366026664
$ws.Range("E19").Value = 0.1677276535728645
$ws.Range("F19").Value = 1.36701075701793
$ws.Range("G19").Value = 0.7606750934688193
$ws.Range("H19").Value = 0.8437285736178382
$ws.Range("I19").Value = 0.9441480764203121
$ws.Range("J19").Value = 0.1828252415163831
$ws.Range("K19").Value = 1.608731644655506
$ws.Range("L19").Value = 0.2168736428828453
$ws.Range("N19").Value = 1.258447817716146
$ws.Range("O19").Value = 3.219122957446586

# Row 20
$ws.Range("C20").Value = 0.2456247615874787
$ws.Range("D20").Value = 0.2058185356969489
$ws.Range("E20").Value = 0.1683202149893397
$ws.Range("F20").Value = 1.366951107918624
$ws.Range("G20").Value = 0.7604122433570382
$ws.Range("H20").Value = 0.8421929325620283
$ws.Range("I20").Value = 0.9439612663641412
$ws.Range("J20").Value = 0.1830316881287928
$ws.Range("K20").Value = 1.658709961526256
$ws.Range("L20").Value = 0.2172580572235248
$ws.Range("N20").Value = 1.261624005677518
$ws.Range("O20").Value = 3.215448485905569

# Row 21
$ws.Range("C21").Value = 0.2512673502673977
$ws.Range("D21").Value = 0.2109159052483278
$ws.Range("E21").Value = 0.1703793770877695
$ws.Range("F21").Value = 1.367577772524058
$ws.Range("G21").Value = 0.7600977338087489
$ws.Range("H21").Value = 0.837459074733232
$ws.Range("I21").Value = 0.9439151788622056
$ws.Range("J21").Value = 0.1838064778195161
$ws.Range("K21").Value = 1.826484405746555
$ws.Range("L21").Value = 0.2186376864015642
$ws.Range("N21").Value = 1.272762968694678
$ws.Range("O21").Value = 3.205187764432907

# Row 22
$ws.Range("C22").Value = 0.2550265609508529
$ws.Range("D22").Value = 0.2142994374299718
$ws.Range("E22").Value = 0.1717757280281305
$ws.Range("F22").Value = 1.368603341175756
$ws.Range("G22").Value = 0.7603163675457694
$ws.Range("H22").Value = 0.8346822642586176
$ws.Range("I22").Value = 0.9443182915925803
$ws.Range("J22").Value = 0.1843731482924014
$ws.Range("K22").Value = 1.936012078227975
$ws.Range("L22").Value = 0.2196047010685334
$ws.Range("N22").Value = 1.280389043453425
$ws.Range("O22").Value = 3.200036114609503

# Row 23
$ws.Range("C23").Value = 0.2530137009568989
$ws.Range("D23").Value = 0.2124888516865298
$ws.Range("E23").Value = 0.1710258606624819
$ws.Range("F23").Value = 1.367999638661686
$ws.Range("G23").Value = 0.760160841836111
$ws.Range("H23").Value = 0.8361353059240741
$ws.Range("I23").Value = 0.9440635159354969
$ws.Range("J23").Value = 0.1840652000971801
$ws.Range("K23").Value = 1.877567210960081
$ws.Range("L23").Value = 0.2190826256705165
$ws.Range("N23").Value = 1.27628733712757
$ws.Range("O23").Value = 3.202643404446178

# Row 24
$ws.Range("C24").Value = 0.2455337514715694
$ws.Range("D24").Value = 0.2057361097160992
$ws.Range("E24").Value = 0.1682874169034321
$ws.Range("F24").Value = 1.366951270991805
$ws.Range("G24").Value = 0.7604245395850597
$ws.Range("H24").Value = 0.8422756605990287
$ws.Range("I24").Value = 0.9439693366485713
$ws.Range("J24").Value = 0.18302004527402
$ws.Range("K24").Value = 1.655965904346431
$ws.Range("L24").Value = 0.2172366150724585
$ws.Range("N24").Value = 1.261447820443223
$ws.Range("O24").Value = 3.215642427380914

# Row 25
$ws.Range("C25").Value = 0.2377761339267721
$ws.Range("D25").Value = 0.1986795840358724
$ws.Range("E25").Value = 0.1655525307653072
$ws.Range("F25").Value = 1.368467878553112
$ws.Range("G25").Value = 0.762527680070832
$ws.Range("H25").Value = 0.8502616905340403
$ws.Range("I25").Value = 0.9457291516360868
$ws.Range("J25").Value = 0.1821527738855053
$ws.Range("K25").Value = 1.416495974944667
$ws.Range("L25").Value = 0.2155277401872056
$ws.Range("N25").Value = 1.246941575512679
$ws.Range("O25").Value = 3.236334940242813
